$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 9.079599999999994
$ws.Range("B6").Value = 6.138900000000003
$ws.Range("B7").Value = 5.585200000000003
$ws.Range("D7").Value = -7.390899999999998
$ws.Range("D12").Value = -7.213400000000002
$ws.Range("E13").Value = 16.7914
$ws.Range("E14").Value = 17.17700000000001
$ws.Range("D15").Value = -8.944099999999992
$ws.Range("B16").Value = 5.0413
$ws.Range("E16").Value = 16.0665
$ws.Range("E19").Value = 16.55880000000001
$ws.Range("B20").Value = 9.314299999999994
$ws.Range("D20").Value = -7.918799999999999
$ws.Range("D21").Value = -8.0871
$ws.Range("D22").Value = -7.5703
$ws.Range("E22").Value = 16.43000000000001
$ws.Range("D23").Value = -7.127399999999995
$ws.Range("B28").Value = 5.977100000000003
$ws.Range("B29").Value = 5.245000000000001
$ws.Range("D29").Value = -7.229699999999995
$ws.Range("B32").Value = 7.291799999999998
$ws.Range("D34").Value = -7.789299999999998
$ws.Range("E36").Value = 15.7739
$ws.Range("B40").Value = 9.183199999999998
$ws.Range("D42").Value = -7.826000000000004
$ws.Range("D43").Value = -8.3148
$ws.Range("D44").Value = -8.131099999999998
$ws.Range("D45").Value = -7.841599999999997
$ws.Range("B46").Value = 6.349400000000004
$ws.Range("D46").Value = -8.411499999999998
$ws.Range("E46").Value = 16.67710000000001
$ws.Range("D50").Value = -8.262699999999997
$ws.Range("E50").Value = 16.8129
$ws.Range("B51").Value = 5.122600000000001
$ws.Range("D51").Value = -7.466899999999999
$ws.Range("B52").Value = 5.014500000000002
$ws.Range("B57").Value = 5.852399999999996
$ws.Range("B59").Value = 5.210299999999998
$ws.Range("B62").Value = 5.596299999999998
$ws.Range("B66").Value = 5.226199999999999
$ws.Range("D66").Value = -7.533200000000001
$ws.Range("D67").Value = -6.631100000000003
$ws.Range("B73").Value = 8.464699999999999
$ws.Range("B74").Value = 9.308399999999988
$ws.Range("D79").Value = -6.228400000000005
$ws.Range("D84").Value = -9.091300000000002
$ws.Range("B92").Value = 4.881099999999996
$ws.Range("D92").Value = -6.592899999999997
$ws.Range("E95").Value = 17.96600000000002
$ws.Range("D97").Value = -8.569899999999997
$ws.Range("E97").Value = 16.65879999999999
$ws.Range("B100").Value = 5.203799999999998
